$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: romiya -------------------------------------------------
$ws.Range("A8:M8").Copy($ws.Range("A9:M9"))

$ws.Range("A9").Value = "romiya"
$ws.Range("B9").Value = 7965458569
$ws.Range("C9").Value = "romiyaji@gmail.com"
$ws.Range("D9").Value = "it"
$ws.Range("E9").Value = "regular"
$ws.Range("F9").Value = 12
$ws.Range("G9").Value = "mahindra"
$ws.Range("H9").Value = 2500000
$ws.Range("I9").Value = 3500000
$ws.Range("J9").Value = "90 days"
$ws.Range("K9").Value = 45728
$ws.Range("L9").Value = "http://example4/resume5.pdf"
$ws.Range("M9").Value = "engineeer"

$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:romiyaji@gmail.com") | Out-Null
$ws.Range("C8").Copy($ws.Range("C9"))
$ws.Range("C9").Value = "romiyaji@gmail.com"

$ws.Hyperlinks.Add($ws.Range("L9"), "http://example4/resume5.pdf") | Out-Null
$ws.Range("L8").Copy($ws.Range("L9"))
$ws.Range("L9").Value = "http://example4/resume5.pdf"

# --- Row 10: somiya -------------------------------------------------
$ws.Range("A8:M8").Copy($ws.Range("A10:M10"))

$ws.Range("A10").Value = "somiya"
$ws.Range("B10").Value = 6965458569
$ws.Range("C10").Value = "somiyaji@gmail.com"
$ws.Range("D10").Value = "it"
$ws.Range("E10").Value = "regular"
$ws.Range("F10").Value = 12
$ws.Range("G10").Value = "mahindra"
$ws.Range("H10").Value = 2500000
$ws.Range("I10").Value = 3500000
$ws.Range("J10").Value = "90 days"
$ws.Range("K10").Value = 45728
$ws.Range("L10").Value = "http://example4/resumesomiya.pdf"
$ws.Range("M10").Value = "engineeer"

$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:somiyaji@gmail.com") | Out-Null
$ws.Range("C8").Copy($ws.Range("C10"))
$ws.Range("C10").Value = "somiyaji@gmail.com"

$ws.Hyperlinks.Add($ws.Range("L10"), "http://example4/resumesomiya.pdf") | Out-Null
$ws.Range("L8").Copy($ws.Range("L10"))
$ws.Range("L10").Value = "http://example4/resumesomiya.pdf"

# --- View cosmetics (best effort) -----------------------------------
$ws.Range("L15").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1

"done"
